$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.174539
$ws.Range("H2").Value = 0.523617
$ws.Range("I2").Value = 0.2427616627057682
$ws.Range("J2").Value = 0.2427616627057681
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.8192526666666667
$ws.Range("N2").Value = 2.457758
$ws.Range("O2").Value = 0.03935738656594369
$ws.Range("P2").Value = 0.03935738656594368
$ws.Range("Q2").Value = 0.1429915411873333
$ws.Range("R2").Value = 1.286923870686
$ws.Range("S2").Value = 0.009554464602502153
$ws.Range("T2").Value = 0.009554464602502149

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.174539
$ws.Range("H3").Value = 0.523617
$ws.Range("I3").Value = 0.2427616627057682
$ws.Range("J3").Value = 0.2427616627057681
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.221608999999999
$ws.Range("N3").Value = 27.664827
$ws.Range("O3").Value = 0.4430115945178313
$ws.Range("P3").Value = 0.4430115945178313
$ws.Range("Q3").Value = 1.609530413251
$ws.Range("R3").Value = 14.485773719259
$ws.Range("S3").Value = 0.1075462312830823
$ws.Range("T3").Value = 0.1075462312830823

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.174539
$ws.Range("H4").Value = 0.523617
$ws.Range("I4").Value = 0.2427616627057682
$ws.Range("J4").Value = 0.2427616627057681
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.77486666666667
$ws.Range("N4").Value = 32.3246
$ws.Range("O4").Value = 0.5176310189162251
$ws.Range("P4").Value = 0.517631018916225
$ws.Range("Q4").Value = 1.880634453133333
$ws.Range("R4").Value = 16.9257100782
$ws.Range("S4").Value = 0.1256609668201837
$ws.Range("T4").Value = 0.1256609668201837

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Agtr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5444336666666666
$ws.Range("H5").Value = 1.633301
$ws.Range("I5").Value = 0.7572383372942318
$ws.Range("J5").Value = 0.7572383372942317
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.8192526666666667
$ws.Range("N5").Value = 2.457758
$ws.Range("O5").Value = 0.03935738656594369
$ws.Range("P5").Value = 0.03935738656594368
$ws.Range("Q5").Value = 0.4460287332397778
$ws.Range("R5").Value = 4.014258599158
$ws.Range("S5").Value = 0.02980292196344154
$ws.Range("T5").Value = 0.02980292196344153

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Agt"
$ws.Range("C6").Value = "Agtr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5444336666666666
$ws.Range("H6").Value = 1.633301
$ws.Range("I6").Value = 0.7572383372942318
$ws.Range("J6").Value = 0.7572383372942317
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.221608999999999
$ws.Range("N6").Value = 27.664827
$ws.Range("O6").Value = 0.4430115945178313
$ws.Range("P6").Value = 0.4430115945178313
$ws.Range("Q6").Value = 5.020554400436332
$ws.Range("R6").Value = 45.18498960392699
$ws.Range("S6").Value = 0.335465363234749
$ws.Range("T6").Value = 0.335465363234749

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Agt"
$ws.Range("C7").Value = "Agtr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5444336666666666
$ws.Range("H7").Value = 1.633301
$ws.Range("I7").Value = 0.7572383372942318
$ws.Range("J7").Value = 0.7572383372942317
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.77486666666667
$ws.Range("N7").Value = 32.3246
$ws.Range("O7").Value = 0.5176310189162251
$ws.Range("P7").Value = 0.517631018916225
$ws.Range("Q7").Value = 5.866200167177777
$ws.Range("R7").Value = 52.7958015046
$ws.Range("S7").Value = 0.3919700520960414
$ws.Range("T7").Value = 0.3919700520960412
